$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.547.07"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.024.95"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.94"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.77"
$ws.Range("E8").Value = "  -8.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.51"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "2.325.18"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.815"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.06"
$ws.Range("E15").Value = "  -7.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.34"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "2.020.08"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "37.520.80"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.55"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "0.0₃0847"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.19"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.04"
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.62"
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.33"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -9.89%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0667"
$ws.Range("E32").Value = "  +6.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.69"
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "1.407.83"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.96"
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.69"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.03"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.31"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "2.217.19"
$ws.Range("E51").Value = "  +1.24%  "

Write-Output "Updated cryptos list"